# Append: 2025-09-28 12:42 JST
# A new scraped listing is inserted as row 3 (pushing the existing rows
# 3-7 down to 4-8), every "取得日時" timestamp in the sheet is refreshed
# to the new run time, and column H is widened to fit the longer skill
# summary text of the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-28 12:42:09"

# --- widen column H (12 -> 18 chars). ColumnWidth adds a fixed
#     0.8333333333333333 padding internally, so back it out to land on
#     an exact stored width of 18. ---
$ws.Columns.Item(8).ColumnWidth = 18 - 0.8333333333333333

# --- insert a new row at position 3; Excel shifts rows 3:7 down to
#     4:8 and carries the hyperlink-style formatting on column F along
#     with it. ---
$ws.Rows.Item(3).Insert()

# --- populate the newly inserted row 3 with the new listing ---
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(3, 2).Value = "【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5394578"
$ws.Cells.Item(3, 7).Value = 158
$ws.Cells.Item(3, 8).Value = "◆自動化,スクレイピング ◇管理"

# --- refresh the "取得日時" timestamp on every data row (2:8) ---
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(5, 1).Value = $newTimestamp
$ws.Cells.Item(6, 1).Value = $newTimestamp
$ws.Cells.Item(7, 1).Value = $newTimestamp
$ws.Cells.Item(8, 1).Value = $newTimestamp

# --- rebuild the F-column hyperlinks in row order (the row insert
#     does not renumber existing hyperlink anchors on its own) ---
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5394578")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5402140")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5402038")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5402182")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5399347")
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5402130")
